# Generate Report for Handoff
# - Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps for the rows that were just (re)handed-off.
# - Marks the "Priority" column as "ht" (handoff type) for the
#   corresponding rows on the zh-cn and de-de localization sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(7, 9, 10, 11, 12, 13)

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $overview.Range("G$r").Value = "2016-09-07 00:33:35"

    # de-de sheet: column H = "Latest Handoff Datetime" (shares the same
    # underlying text as the Overview sheet's column G)
    $dede.Range("H$r").Value = "2016-09-07 00:33:35"

    # zh-cn sheet: column H = "Latest Handoff Datetime"
    $zhcn.Range("H$r").Value = "2016-09-07 00:33:30"

    # zh-cn / de-de sheets: column E = "Priority" -> "ht"
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
